# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header/style conventions and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, using the same style as the other header cells (s="1").
# Copy format from the existing "IP" header (H1) so the bold/border/
# alignment formatting carries over to the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row values for columns I and J (rows 2..67)
$iValues = @(7,9,9,8,6,7,8,7,8,7,7,7,10,7,6,8,6,9,8,8,7,7,8,7,6,7,7,7,8,7,6,6,7,6,8,8,7,8,6,7,8,7,8,7,6,6,6,8,6,7,7,10,9,6,6,8,7,7,7,7,7,8,8,7,8,7)
$jValues = @(7,9,9,8,7,7,8,8,8,7,7,7,10,7,7,8,7,9,8,8,7,8,8,7,7,7,7,8,8,7,6,6,7,6,8,8,8,8,6,7,8,7,8,7,6,6,6,8,6,7,7,10,9,7,7,8,7,7,7,8,7,8,8,7,8,7)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
